$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 871.25
$ws.Range("J8").Value = 999.3333
$ws.Range("L8").Value = 2997.9999
$ws.Range("N8").Value = -3275.9999
$ws.Range("H17").Value = 2188
$ws.Range("J17").Value = 2188
$ws.Range("L17").Value = 6564
$ws.Range("N17").Value = -6900
$ws.Range("H32").Value = 3466.6667
$ws.Range("I32").Value = 3450.5
$ws.Range("K32").Value = 3450.5
$ws.Range("M32").Value = -3124.5
$ws.Range("H39").Value = 2452.9167
$ws.Range("I39").Value = 998.25
$ws.Range("J39").Value = 5362.25
$ws.Range("K39").Value = 2994.75
$ws.Range("L39").Value = 16086.75
$ws.Range("M39").Value = -2698.75
$ws.Range("N39").Value = -16678.75
$ws.Range("H40").Value = 2615.6667
$ws.Range("I40").Value = 1350
$ws.Range("K40").Value = 1350
$ws.Range("M40").Value = -1175
$ws.Range("H48").Value = 4586.25
$ws.Range("J48").Value = 4586.25
$ws.Range("L48").Value = 13758.75
$ws.Range("N48").Value = -14342.75
$ws.Range("H56").Value = 4586.25
$ws.Range("J56").Value = 4586.25
$ws.Range("L56").Value = 13758.75
$ws.Range("N56").Value = -14826.75
$ws.Range("H62").Value = 1647.3334
$ws.Range("I62").Value = 1647.3334
$ws.Range("K62").Value = 1647.3334
$ws.Range("M62").Value = -1023.3334
$ws.Range("H65").Value = 1647.3334
$ws.Range("I65").Value = 1647.3334
$ws.Range("K65").Value = 8236.666999999999
$ws.Range("M65").Value = -5116.666999999999
$ws.Range("H70").Value = 11699.363
$ws.Range("I70").Value = 3099.5
$ws.Range("J70").Value = 16613.572
$ws.Range("K70").Value = 9298.5
$ws.Range("L70").Value = 49840.716
$ws.Range("M70").Value = -9028.5
$ws.Range("N70").Value = -50380.716
$ws.Range("H73").Value = 11699.363
$ws.Range("I73").Value = 3099.5
$ws.Range("J73").Value = 16613.572
$ws.Range("K73").Value = 9298.5
$ws.Range("L73").Value = 49840.716
$ws.Range("M73").Value = -8362.5
$ws.Range("N73").Value = -51712.716
$ws.Range("H96").Value = 638.3333
$ws.Range("I96").Value = 536
$ws.Range("J96").Value = 996.5
$ws.Range("K96").Value = 1608
$ws.Range("L96").Value = 2989.5
$ws.Range("M96").Value = -235
$ws.Range("N96").Value = -5735.5
$ws.Range("H98").Value = 9323.166999999999
$ws.Range("I98").Value = 897.625
$ws.Range("J98").Value = 26174.25
$ws.Range("K98").Value = 897.625
$ws.Range("L98").Value = 26174.25
$ws.Range("M98").Value = 600.375
$ws.Range("N98").Value = -29170.25
$ws.Range("H99").Value = 2061.8572
$ws.Range("I99").Value = 483.5
$ws.Range("K99").Value = 1450.5
$ws.Range("M99").Value = 47.5
$ws.Range("H100").Value = 10350.7
$ws.Range("I100").Value = 2492.5
$ws.Range("J100").Value = 12315.25
$ws.Range("K100").Value = 2492.5
$ws.Range("L100").Value = 12315.25
$ws.Range("M100").Value = -1951.5
$ws.Range("N100").Value = -13397.25
$ws.Range("H101").Value = 2939.4443
$ws.Range("I101").Value = 1126
$ws.Range("J101").Value = 6566.3335
$ws.Range("K101").Value = 3378
$ws.Range("L101").Value = 19699.0005
$ws.Range("M101").Value = -1756
$ws.Range("N101").Value = -22943.0005
$ws.Range("H111").Value = 2928.5
$ws.Range("I111").Value = 2786
$ws.Range("J111").Value = 3498.5
$ws.Range("K111").Value = 8358
$ws.Range("L111").Value = 10495.5
$ws.Range("M111").Value = -5291
$ws.Range("N111").Value = -16629.5
$ws.Range("H122").Value = 9323.166999999999
$ws.Range("I122").Value = 897.625
$ws.Range("J122").Value = 26174.25
$ws.Range("K122").Value = 2692.875
$ws.Range("L122").Value = 78522.75
$ws.Range("M122").Value = -242.875
$ws.Range("N122").Value = -83422.75
$ws.Range("H125").Value = 3921.4375
$ws.Range("I125").Value = 3794.4443
$ws.Range("J125").Value = 4084.7144
$ws.Range("K125").Value = 34149.9987
$ws.Range("L125").Value = 36762.4296
$ws.Range("M125").Value = -31689.9987
$ws.Range("N125").Value = -41682.4296
$ws.Range("H129").Value = 1618.6923
$ws.Range("J129").Value = 3659.3333
$ws.Range("L129").Value = 10977.9999
$ws.Range("N129").Value = -20977.9999
$ws.Range("H137").Value = 991.1
$ws.Range("I137").Value = 285.5
$ws.Range("J137").Value = 1461.5
$ws.Range("K137").Value = 856.5
$ws.Range("L137").Value = 4384.5
$ws.Range("M137").Value = 1693.5
$ws.Range("N137").Value = -9484.5
$ws.Range("H138").Value = 5510.1934
$ws.Range("I138").Value = 6316.923
$ws.Range("J138").Value = 4927.5557
$ws.Range("K138").Value = 18950.769
$ws.Range("L138").Value = 14782.6671
$ws.Range("M138").Value = -13810.769
$ws.Range("N138").Value = -25062.6671
$ws.Range("H141").Value = 14010.4
$ws.Range("I141").Value = 6649.3335
$ws.Range("J141").Value = 25052
$ws.Range("K141").Value = 19948.0005
$ws.Range("L141").Value = 75156
$ws.Range("M141").Value = -14768.0005
$ws.Range("N141").Value = -85516

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2840.9429
$ws.Range("I32").Value = 2013.4849
$ws.Range("K32").Value = 2013.4849
$ws.Range("M32").Value = -1726.4849
$ws.Range("H45").Value = 3286.4167
$ws.Range("I45").Value = 3005
$ws.Range("K45").Value = 3005
$ws.Range("M45").Value = -2628
$ws.Range("H61").Value = 5253.8
$ws.Range("I61").Value = 4650.75
$ws.Range("J61").Value = 7666
$ws.Range("K61").Value = 4650.75
$ws.Range("L61").Value = 7666
$ws.Range("M61").Value = -4438.75
$ws.Range("N61").Value = -8090
$ws.Range("H74").Value = 2500.4348
$ws.Range("I74").Value = 1636.5834
$ws.Range("K74").Value = 1636.5834
$ws.Range("M74").Value = -762.5834
$ws.Range("H77").Value = 2500.4348
$ws.Range("I77").Value = 1636.5834
$ws.Range("K77").Value = 8182.916999999999
$ws.Range("M77").Value = -3814.916999999999
$ws.Range("H102").Value = 3983.8
$ws.Range("I102").Value = 3977
$ws.Range("K102").Value = 3977
$ws.Range("M102").Value = -2355
$ws.Range("H110").Value = 1398.7646
$ws.Range("I110").Value = 1361.1875
$ws.Range("K110").Value = 1361.1875
$ws.Range("M110").Value = 683.8125
$ws.Range("H122").Value = 1766.5
$ws.Range("I122").Value = 533.3333
$ws.Range("J122").Value = 2999.6667
$ws.Range("K122").Value = 1599.9999
$ws.Range("L122").Value = 8999.000100000001
$ws.Range("M122").Value = 850.0001
$ws.Range("N122").Value = -13899.0001
$ws.Range("H132").Value = 1517.2106
$ws.Range("I132").Value = 1393.1945
$ws.Range("J132").Value = 3749.5
$ws.Range("K132").Value = 4179.583500000001
$ws.Range("L132").Value = 11248.5
$ws.Range("M132").Value = -1649.583500000001
$ws.Range("N132").Value = -16308.5
$ws.Range("H136").Value = 5253.8
$ws.Range("I136").Value = 4650.75
$ws.Range("J136").Value = 7666
$ws.Range("K136").Value = 13952.25
$ws.Range("L136").Value = 22998
$ws.Range("M136").Value = -11402.25
$ws.Range("N136").Value = -28098

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 450
$ws.Range("I22").Value = 450
$ws.Range("K22").Value = 450
$ws.Range("M22").Value = -277
$ws.Range("H94").Value = 19776.637
$ws.Range("I94").Value = 14923
$ws.Range("J94").Value = 25601
$ws.Range("K94").Value = 14923
$ws.Range("L94").Value = 25601
$ws.Range("M94").Value = -14472
$ws.Range("N94").Value = -26503
$ws.Range("H105").Value = 3796.182
$ws.Range("I105").Value = 3208.5715
$ws.Range("K105").Value = 3208.5715
$ws.Range("M105").Value = -1461.5715
$ws.Range("H107").Value = 1047.6666
$ws.Range("I107").Value = 996.8889
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 996.8889
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 923.1111
$ws.Range("N107").Value = -5040
$ws.Range("H134").Value = 8429.117
$ws.Range("I134").Value = 8556.286
$ws.Range("J134").Value = 7835.6665
$ws.Range("K134").Value = 25668.858
$ws.Range("L134").Value = 23506.9995
$ws.Range("M134").Value = -23133.858
$ws.Range("N134").Value = -28576.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3963.9375
$ws.Range("I31").Value = 3096.9285
$ws.Range("J31").Value = 4638.278
$ws.Range("K31").Value = 3096.9285
$ws.Range("L31").Value = 4638.278
$ws.Range("M31").Value = -2801.9285
$ws.Range("N31").Value = -5228.278
$ws.Range("H34").Value = 3963.9375
$ws.Range("I34").Value = 3096.9285
$ws.Range("J34").Value = 4638.278
$ws.Range("K34").Value = 3096.9285
$ws.Range("L34").Value = 4638.278
$ws.Range("M34").Value = -2894.9285
$ws.Range("N34").Value = -5042.278
$ws.Range("H58").Value = 2548.0833
$ws.Range("I58").Value = 1451.8572
$ws.Range("K58").Value = 1451.8572
$ws.Range("M58").Value = -1248.8572
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("H86").Value = 6866.421
$ws.Range("I86").Value = 6480.25
$ws.Range("J86").Value = 7147.273
$ws.Range("K86").Value = 6480.25
$ws.Range("L86").Value = 7147.273
$ws.Range("M86").Value = -5357.25
$ws.Range("N86").Value = -9393.273000000001
$ws.Range("H89").Value = 6866.421
$ws.Range("I89").Value = 6480.25
$ws.Range("J89").Value = 7147.273
$ws.Range("K89").Value = 32401.25
$ws.Range("L89").Value = 35736.365
$ws.Range("M89").Value = -26785.25
$ws.Range("N89").Value = -46968.365
$ws.Range("H122").Value = 5126
$ws.Range("I122").Value = 3501.3333
$ws.Range("K122").Value = 10503.9999
$ws.Range("M122").Value = -8053.999899999999
$ws.Range("H134").Value = 2449.5
$ws.Range("I134").Value = 2449.5
$ws.Range("K134").Value = 7348.5
$ws.Range("M134").Value = -4813.5
$ws.Range("H136").Value = 2548.0833
$ws.Range("I136").Value = 1451.8572
$ws.Range("K136").Value = 4355.571599999999
$ws.Range("M136").Value = -1805.571599999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 63245370
$ws.Range("J4").Value = 1000000
$ws.Range("L4").Value = 3000000
$ws.Range("N4").Value = -3000224
$ws.Range("H36").Value = 500
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 500
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = 1500
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -1838
$ws.Range("H44").Value = 391.66666
$ws.Range("I44").Value = 391.66666
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 1174.99998
$ws.Range("L44").Value = 0
$ws.Range("M44").Value = -776.9999800000001
$ws.Range("N44").ClearContents()
$ws.Range("H57").Value = 8899.200000000001
$ws.Range("J57").Value = 9927.429
$ws.Range("L57").Value = 29782.287
$ws.Range("N57").Value = -30900.287
$ws.Range("H68").Value = 1211
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1211
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3633
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -5255
$ws.Range("H71").Value = 1211
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1211
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 10899
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -19011
$ws.Range("H80").Value = 3499.6
$ws.Range("H83").Value = 3499.6
$ws.Range("H103").Value = 9773.817999999999
$ws.Range("I103").Value = 13254.25
$ws.Range("J103").Value = 492.66666
$ws.Range("K103").Value = 39762.75
$ws.Range("L103").Value = 1477.99998
$ws.Range("M103").Value = -38883.75
$ws.Range("N103").Value = -3235.99998
$ws.Range("H114").Value = 2043.5
$ws.Range("I114").Value = 158.5
$ws.Range("K114").Value = 475.5
$ws.Range("M114").Value = 2778.5
$ws.Range("H129").Value = 2812.4062
$ws.Range("J129").Value = 3006.7932
$ws.Range("L129").Value = 9020.3796
$ws.Range("N129").Value = -19020.3796
$ws.Range("H131").Value = 1579.8695
$ws.Range("J131").Value = 2189.6365
$ws.Range("L131").Value = 6568.9095
$ws.Range("N131").Value = -16648.9095

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 9111.294
$ws.Range("I102").Value = 10366.296
$ws.Range("K102").Value = 10366.296
$ws.Range("M102").Value = -8744.296
$ws.Range("H113").Value = 4895.769
$ws.Range("I113").Value = 3377.4
$ws.Range("K113").Value = 3377.4
$ws.Range("M113").Value = -1207.4
$ws.Range("H122").Value = 4152.7
$ws.Range("I122").Value = 3508.0557
$ws.Range("J122").Value = 5119.6665
$ws.Range("K122").Value = 10524.1671
$ws.Range("L122").Value = 15358.9995
$ws.Range("M122").Value = -8074.167099999999
$ws.Range("N122").Value = -20258.9995
$ws.Range("H126").Value = 8754
$ws.Range("I126").Value = 7827.143
$ws.Range("J126").Value = 10916.667
$ws.Range("K126").Value = 23481.429
$ws.Range("L126").Value = 32750.001
$ws.Range("M126").Value = -21011.429
$ws.Range("N126").Value = -37690.001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3433.25
$ws.Range("I7").Value = 3736
$ws.Range("J7").Value = 2525
$ws.Range("K7").Value = 3736
$ws.Range("L7").Value = 2525
$ws.Range("M7").Value = -3624
$ws.Range("N7").Value = -2749
$ws.Range("H22").Value = 1414.3704
$ws.Range("I22").Value = 1128
$ws.Range("J22").Value = 1479.4546
$ws.Range("K22").Value = 1128
$ws.Range("L22").Value = 1479.4546
$ws.Range("M22").Value = -833
$ws.Range("N22").Value = -2069.4546
$ws.Range("H27").Value = 1414.3704
$ws.Range("I27").Value = 1128
$ws.Range("J27").Value = 1479.4546
$ws.Range("K27").Value = 1128
$ws.Range("L27").Value = 1479.4546
$ws.Range("M27").Value = -1021
$ws.Range("N27").Value = -1693.4546
$ws.Range("H40").Value = 2533.4
$ws.Range("I40").Value = 2164.5
$ws.Range("K40").Value = 2164.5
$ws.Range("M40").Value = -2028.5
$ws.Range("H46").Value = 1543.2858
$ws.Range("I46").Value = 1079.6
$ws.Range("J46").Value = 1688.1875
$ws.Range("K46").Value = 1079.6
$ws.Range("L46").Value = 1688.1875
$ws.Range("M46").Value = -891.5999999999999
$ws.Range("N46").Value = -2064.1875
$ws.Range("H61").Value = 3558.7273
$ws.Range("I61").Value = 3564.7
$ws.Range("K61").Value = 3564.7
$ws.Range("M61").Value = -3362.7
$ws.Range("H108").Value = 64999.5
$ws.Range("J108").Value = 64999.5
$ws.Range("L108").Value = 64999.5
$ws.Range("N108").Value = -72679.5
$ws.Range("H113").Value = 3558.7273
$ws.Range("I113").Value = 3564.7
$ws.Range("K113").Value = 3564.7
$ws.Range("M113").Value = -1394.7
$ws.Range("H122").Value = 6865.2
$ws.Range("I122").Value = 6920.25
$ws.Range("J122").Value = 6828.5
$ws.Range("K122").Value = 20760.75
$ws.Range("L122").Value = 20485.5
$ws.Range("M122").Value = -18310.75
$ws.Range("N122").Value = -25385.5
$ws.Range("H126").Value = 3433.25
$ws.Range("I126").Value = 3736
$ws.Range("J126").Value = 2525
$ws.Range("K126").Value = 11208
$ws.Range("L126").Value = 7575
$ws.Range("M126").Value = -8738
$ws.Range("N126").Value = -12515
$ws.Range("H132").Value = 2637.7368
$ws.Range("I132").Value = 2306.6667
$ws.Range("J132").Value = 2790.5386
$ws.Range("K132").Value = 6920.000100000001
$ws.Range("L132").Value = 8371.6158
$ws.Range("M132").Value = -4390.000100000001
$ws.Range("N132").Value = -13431.6158
$ws.Range("H136").Value = 2252.6
$ws.Range("I136").Value = 2123.25
$ws.Range("J136").Value = 2446.625
$ws.Range("K136").Value = 6369.75
$ws.Range("L136").Value = 7339.875
$ws.Range("M136").Value = -3819.75
$ws.Range("N136").Value = -12439.875

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 864.6667
$ws.Range("I100").Value = 856.25
$ws.Range("J100").Value = 932
$ws.Range("K100").Value = 1712.5
$ws.Range("L100").Value = 1864
$ws.Range("M100").Value = -1171.5
$ws.Range("N100").Value = -2946
$ws.Range("H126").Value = 1982.0555
$ws.Range("I126").Value = 1888.6666
$ws.Range("J126").Value = 2168.8333
$ws.Range("K126").Value = 5665.9998
$ws.Range("L126").Value = 6506.499899999999
$ws.Range("M126").Value = -3195.9998
$ws.Range("N126").Value = -11446.4999
$ws.Range("H136").Value = 2591.8235
$ws.Range("I136").Value = 2029.1154
$ws.Range("K136").Value = 6087.3462
$ws.Range("M136").Value = -3537.3462
